# Updated cryptos list with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price ("D") cells hold text that looks numeric (dotted thousands, e.g. "27.501.81").
# Force them through the Text number format on write so Excel COM does not
# auto-coerce the assignment into a real number, then restore the default style
# (matches the original cells, which carry no explicit style index).

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.501.81"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -1.53%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.842.05"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -2.01%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.005"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -1.38%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "333.23"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.59%  "
$ws.Range("E6").Value = "  -1.20%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4638"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.74%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3846"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.58%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "45.89"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.88%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07879"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.72%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.9934"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.13%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "21.44"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.54%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.851.86"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -2.64%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.934"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.11%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.104"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.24%  "
$ws.Range("E16").Value = "  -1.35%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "88.59"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.32%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.06665"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.93%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.00001034"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.94%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.08"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.64%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.003"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.37%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "27.533.44"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.48%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.372"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.59%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "10.88"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.63%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.303"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -2.33%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.074.69"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -2.25%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "158.52"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.60%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "19.44"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -2.29%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.109"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +2.07%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "5.391"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.01%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "119.57"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.87%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.9732"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +2.13%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.09388"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.39%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.591"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.93%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.287"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.48%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.334"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.74%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.06022"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.44%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02223"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.49%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "8.299"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +2.63%  "
$ws.Range("E40").Value = "  -1.81%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.5875"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.15%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1863"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.51%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "10.32"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.86%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.233"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -3.08%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.5579"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.76%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "12.18"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.94%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.902"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.46%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.06692"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -2.40%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "110.82"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -2.27%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.049"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.39%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.004"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.47%  "
